$wb = $excel.ActiveWorkbook

# -----------------------------------------------------------------
# Sheet "plane": append a new last column EE ("h_is") after ED.
# -----------------------------------------------------------------
$wsPlane = $wb.Worksheets.Item("plane")

# Copy the formatting from the previous last column (ED) into the new EE
# column so the header keeps the bold/border/centered style.
$wsPlane.Range("ED1:ED5").Copy()
$wsPlane.Range("EE1:EE5").PasteSpecial(-4122)

$wsPlane.Range("EE1").Value = "h_is"
$wsPlane.Range("EE2").Value = "b''"
$wsPlane.Range("EE3").Value = 490189.7099025669
# EE4 stays blank (mirrors ED4, an empty placeholder cell for this row);
# the PasteSpecial above already materialized the cell with ED4's format,
# so no explicit value assignment is needed here.
$wsPlane.Range("EE5").Value = 474663.510268038

# -----------------------------------------------------------------
# Sheet "cascade": insert a new "h_is_throat" column right before the
# existing "dh_s" column, and drop the five trailing
# "efficiency_drop_*" columns.
# -----------------------------------------------------------------
$wsCascade = $wb.Worksheets.Item("cascade")

# Insert a blank column before EL ("dh_s"); this pushes
# dh_s, incidence, efficiency_drop_* one column to the right
# (EL->EM, EM->EN, ..., ER->ES).
$wsCascade.Range("EL1").EntireColumn.Insert()

# Match formatting of the neighboring "blockage_throat" column (EK).
$wsCascade.Range("EK1:EK3").Copy()
$wsCascade.Range("EL1:EL3").PasteSpecial(-4122)

$wsCascade.Range("EL1").Value = "h_is_throat"
$wsCascade.Range("EL2").Value = 490646.156185923
$wsCascade.Range("EL3").Value = 478710.3421281012

# Remove the five "efficiency_drop_*" columns, which after the insert
# above now sit at EO:ES.
$wsCascade.Range("EO1:ES1").EntireColumn.Delete()

# -----------------------------------------------------------------
# Sheet "geometry": column "gauging_angle" (L) moves from the front of
# the cascade-geometry block to the very end (after
# "leading_edge_diameter_chord_ratio"), i.e. the L:AU block is
# rotated one column to the left.
# -----------------------------------------------------------------
$wsGeom = $wb.Worksheets.Item("geometry")

# Remember the exact values of column L (gauging_angle) before it is
# removed - read them now so they can be written back after the shift.
$gaugingHeader = $wsGeom.Range("L1").Value2
$gaugingRow2 = $wsGeom.Range("L2").Value2
$gaugingRow3 = $wsGeom.Range("L3").Value2

# Copy column L's formatting (bold header / plain data style) so it can
# be reapplied to the new last column after the shift.
$wsGeom.Range("L1:L3").Copy()

# Delete column L; M:AU shifts left into L:AT.
$wsGeom.Range("L1").EntireColumn.Delete()

# Paste the remembered formatting into the new last column (AU), then
# write the "gauging_angle" values back into it.
$wsGeom.Range("AU1:AU3").PasteSpecial(-4122)
$wsGeom.Range("AU1").Value = $gaugingHeader
$wsGeom.Range("AU2").Value = $gaugingRow2
$wsGeom.Range("AU3").Value = $gaugingRow3
